$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row's shared-string labels:
#    "..._old" -> "..._FV2310" and "..._new" -> "..._FV2404"
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $ws.Range("$($leftCols[$i])1").Value  = "$($oldHeaders[$i])_FV2310"
    $ws.Range("$($rightCols[$i])1").Value = "$($oldHeaders[$i])_FV2404"
}

# 2. Turn the data range into an Excel Table ("Table1") spanning A1:U64
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3. Freeze the header row (pane split after row 1)
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
